# Update the "2024" sheet: a new transaction row is inserted above the
# existing row 10 ("Loan" group, August "amazeloan" entry), pushing the
# rest of the Loan sub-table (rows 10-32) down by one row (to 11-33).
# The newly inserted row 10 gets a brand-new September "amazeloan" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at position 10; everything from the old row 10
# down to row 32 shifts down to 11-33 (dimension grows from Y32 to Y33).
$ws.Rows("10").Insert()

# Populate the newly-inserted row 10 with the new September entry.
$ws.Range("R10").Value = "amazeloan"
$ws.Range("S10").Value = "2024-09-01 10:27:41"
